$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 157; this shifts the existing rows 157-229
# down to 158-230, preserving all their data untouched.
$ws.Rows("157:157").Insert()

# Populate the newly inserted row 157 with the new record.
$ws.Range("A157").Value = 3
$ws.Range("B157").Value = "Femacal de La Calera"
$ws.Range("C157").Value = "Coquimbo"
$ws.Range("D157").Value = 44510
$ws.Range("E157").Value = 5
$ws.Range("F157").Value = 100112043
$ws.Range("G157").Value = "Pepino ensalada"
$ws.Range("H157").Value = "Sin especificar"
$ws.Range("I157").Value = "Primera"
$ws.Range("J157").Value = 130
$ws.Range("K157").Value = 7000
$ws.Range("L157").Value = 7500
$ws.Range("M157").Value = 7231
$ws.Range("N157").Value = "$/caja 70 unidades"
$ws.Range("O157").Value = "Región de Arica y Parinacota"
$ws.Range("P157").Value = 103
$ws.Range("Q157").Value = 70
$ws.Range("R157").Value = "Hortaliza"
